$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "ceasemultipleserviceaccount"
$ws.Range("B21").Value = "280002720000"

$ws.Range("A22").Value = "CountofServiceProducts"
$ws.Range("B22").Value = "9177989"

$ws.Range("F17").Select()

$ws.Columns.Item(1).ColumnWidth = 26.5
